$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44239
$ws.Range("K2").Value = "Fortuna"
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S2").Value = 861
$ws.Range("D3").Value = 44169
$ws.Range("K3").Value = "Angeleno"
$ws.Range("L3").Value = "Tercera"
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("S3").Value = 1361
$ws.Range("D4").Value = 44243
$ws.Range("K4").Value = "Black Amber"
$ws.Range("L4").Value = "Primera"
$ws.Range("Q4").Value = "`$/caja 18 kilos granel"
$ws.Range("D5").Value = 44580
$ws.Range("K5").Value = "Black Amber"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1083
$ws.Range("D6").Value = 44587
$ws.Range("K6").Value = "Black Amber"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "`$/caja 18 kilos granel"
$ws.Range("S6").Value = 861
$ws.Range("D7").Value = 44278
$ws.Range("K7").Value = "Angeleno"
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("Q7").Value = "`$/caja 18 kilos granel"
$ws.Range("S7").Value = 861
$ws.Range("D8").Value = 44314
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S8").Value = 806
$ws.Range("D9").Value = 44245
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S9").Value = 806
$ws.Range("D10").Value = 44574
$ws.Range("L10").Value = "Primera"
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 18500
$ws.Range("S10").Value = 1028
$ws.Range("D11").Value = 44217
$ws.Range("K11").Value = "Black Amber"
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 917
$ws.Range("D12").Value = 44238
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14500
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 806
$ws.Range("D13").Value = 44238
$ws.Range("K13").Value = "Fortuna"
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 14500
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 806
$ws.Range("D15").Value = 44174
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20500
$ws.Range("Q15").Value = "`$/caja 18 kilos granel"
$ws.Range("S15").Value = 1139
$ws.Range("D16").Value = 44229
$ws.Range("K16").Value = "Fortuna"
$ws.Range("L16").Value = "Segunda"
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = "`$/bandeja 18 kilos granel"
$ws.Range("S16").Value = 806
$ws.Range("D17").Value = 44285
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 300
